# Append the new data row (row 74) to the profit-data worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 74

# Column A holds the date as literal text (matches existing rows, which are
# inline/shared strings like "02/05/2026", not real date serials), so force
# a text format before assigning to stop Excel auto-converting it to a date.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "02/06/2026"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 8731.5
$ws.Cells.Item($row, 3).Value = 0.2450454044534087
$ws.Cells.Item($row, 4).Value = 0.7549545955465913
$ws.Cells.Item($row, 5).Value = -353.76
$ws.Cells.Item($row, 6).Value = -42.37
$ws.Cells.Item($row, 7).Value = -24299.87
$ws.Cells.Item($row, 8).Value = -78.66
$ws.Cells.Item($row, 9).Value = -1163.24
$ws.Cells.Item($row, 10).Value = -35.22
$ws.Cells.Item($row, 11).Value = -25463.11
$ws.Cells.Item($row, 12).Value = -74.47
